$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) store text-like values (e.g. "28.242.05",
# "  +2.72%  ") as inline strings, not numbers. Force the target ranges to a
# text number format before assigning the new values so that Excel keeps
# storing them as text instead of silently reinterpreting/rounding them as
# numeric values (which would corrupt values like "0.000008867" or lose
# trailing zeros like "318.08").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.242.05'
$ws.Range("E2").Value = '  +2.72%  '
$ws.Range("D3").Value = '1.918.98'
$ws.Range("E3").Value = '  +2.50%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.97%  '
$ws.Range("D5").Value = '318.08'
$ws.Range("E5").Value = '  +1.42%  '
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("D7").Value = '0.4847'
$ws.Range("E7").Value = '  +1.12%  '
$ws.Range("D8").Value = '0.3842'
$ws.Range("E8").Value = '  +2.01%  '
$ws.Range("D9").Value = '0.07396'
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").Value = '0.9430'
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("D11").Value = '20.99'
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("D12").Value = '0.07822'
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").Value = '1.929.80'
$ws.Range("E13").Value = '  +2.65%  '
$ws.Range("D14").Value = '5.522'
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").Value = '6.661'
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("D16").Value = '91.55'
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("D17").Value = '1.007'
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").Value = '0.000008867'
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").Value = '1.005'
$ws.Range("E19").Value = '  -0.92%  '
$ws.Range("D20").Value = '28.255.16'
$ws.Range("E20").Value = '  +2.66%  '
$ws.Range("D21").Value = '14.91'
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").Value = '5.174'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").Value = '2.177.25'
$ws.Range("E23").Value = '  +2.75%  '
$ws.Range("D24").Value = '10.95'
$ws.Range("E24").Value = '  +2.30%  '
$ws.Range("D25").Value = '156.47'
$ws.Range("E25").Value = '  +1.57%  '
$ws.Range("D26").Value = '1.926'
$ws.Range("E26").Value = '  -1.92%  '
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("D28").Value = '2.111'
$ws.Range("E28").Value = '  +4.57%  '
$ws.Range("D29").Value = '116.78'
$ws.Range("D30").Value = '4.999'
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("D31").Value = '0.08927'
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Value = '3.360'
$ws.Range("E32").Value = '  +0.95%  '
$ws.Range("E33").Value = '  +3.04%  '
$ws.Range("D34").Value = '0.7781'
$ws.Range("E34").Value = '  +3.92%  '
$ws.Range("D35").Value = '4.714'
$ws.Range("E35").Value = '  +2.44%  '
$ws.Range("D36").Value = '2.701'
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("D37").Value = '0.02056'
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("D38").Value = '1.107'
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("D39").Value = '0.5566'
$ws.Range("E39").Value = '  +3.50%  '
$ws.Range("D40").Value = '0.05342'
$ws.Range("E40").Value = '  +0.86%  '
$ws.Range("D41").Value = '3.010'
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("D42").Value = '7.066'
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("D43").Value = '0.1536'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").Value = '8.517'
$ws.Range("E44").Value = '  +0.99%  '
$ws.Range("D45").Value = '10.81'
$ws.Range("E45").Value = '  +1.99%  '
$ws.Range("D46").Value = '0.4887'
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("D47").Value = '107.22'
$ws.Range("E47").Value = '  +3.89%  '
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("D49").Value = '1.666'
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").Value = '68.88'
$ws.Range("E50").Value = '  +2.59%  '
$ws.Range("D51").Value = '0.06125'
$ws.Range("E51").Value = '  +0.18%  '
